$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the location / id data between rows 3 and 4
$a3 = $ws.Range("A3").Value2
$q3 = $ws.Range("Q3").Value2
$r3 = $ws.Range("R3").Value2

$a4 = $ws.Range("A4").Value2
$q4 = $ws.Range("Q4").Value2
$r4 = $ws.Range("R4").Value2

$ws.Range("A3").Value2 = $a4
$ws.Range("Q3").Value2 = $q4
$ws.Range("R3").Value2 = $r4

$ws.Range("A4").Value2 = $a3
$ws.Range("Q4").Value2 = $q3
$ws.Range("R4").Value2 = $r3

# Swap the location / id / comment data between rows 13 and 15
$a13 = $ws.Range("A13").Value2
$q13 = $ws.Range("Q13").Value2
$r13 = $ws.Range("R13").Value2
$ac13 = $ws.Range("AC13").Value2

$a15 = $ws.Range("A15").Value2
$q15 = $ws.Range("Q15").Value2
$r15 = $ws.Range("R15").Value2
$ac15 = $ws.Range("AC15").Value2

$ws.Range("A13").Value2 = $a15
$ws.Range("Q13").Value2 = $q15
$ws.Range("R13").Value2 = $r15
$ws.Range("AC13").Value2 = $ac15

$ws.Range("A15").Value2 = $a13
$ws.Range("Q15").Value2 = $q13
$ws.Range("R15").Value2 = $r13
$ws.Range("AC15").Value2 = $ac13
